$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New participants that failed QC manually (all FALSE) plus one new entry
# (sub_025) that failed with a reason noted in column C.
$newRows = @(
    @("sub_021", $false, $null),
    @("sub_022", $false, $null),
    @("sub_023", $false, $null),
    @("sub_024", $false, $null),
    @("sub_025", $true,  "PC froze")
)

$startRow = 22
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    if ($row[2] -ne $null) {
        $ws.Cells.Item($r, 3).Value = $row[2]
    }
}

$ws.Range("C27").Select()
